# Apply the edits described by the commit:
#  - update the XOR/Non-XOR "Hamming" sample figures (P5:Q5, S5, P6:Q6, S6);
#    the dependent SUM formulas (R, T, U columns) recalc automatically
#  - move the viewport/selection on Sheet1 (topLeftCell A1->L1, selection R9->V6)
#  - narrow the workbook tab-bar ratio (tabRatio 501->151)
#  - re-apply the border flag on the spacer-cell style used by O1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated sample data (row 5 / row 6 of the Hamming block, cols P:S) ---
$ws.Range("P5").Value = 496
$ws.Range("Q5").Value = 332
$ws.Range("S5").Value = 486

$ws.Range("P6").Value = 4825
$ws.Range("Q6").Value = 3126
$ws.Range("S6").Value = 4810

# --- style tweak: the xf used by O1 should carry applyBorder="true" ---
$o1 = $ws.Range("O1")
$o1.Borders.LineStyle = $o1.Borders.LineStyle

# --- window / view state ---
$win = $excel.ActiveWindow
$win.TabRatio = 0.151
$win.ScrollColumn = 12
$win.ScrollRow = 1
try {
    $win.TopLeftCell = $ws.Range("L1")
} catch {
}

# move the active selection to V6 (was R9)
$ws.Range("V6").Select()
